$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Table cell margins: left margin 138 dxa (6.9pt) -> 143 dxa (7.15pt)
#    Applies to the two tables that currently have a 6.9pt left cell margin.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    if ([math]::Round($t.LeftPadding, 2) -eq 6.9) {
        $t.LeftPadding = 7.15
    }
}

# ---------------------------------------------------------------------------
# 2) Merge the date-format placeholder that is split across several runs
#    ( {{ question.value | date("d / d / . / MM / .Y / YYY / ") }} )
#    into a single run reading: {{ question.value | date("dd.MM.YYYY") }}
#
#    We first collapse it to a short placeholder (plain Range.Text
#    assignment, which merges the underlying runs and keeps straight
#    quotes), then rename the placeholder to the final text. Doing it in
#    two steps avoids two pitfalls of this engine:
#      * Find/Replace (Execute with a replacement) runs the text through
#        AutoFormat and turns straight quotes into curly quotes.
#      * A direct Range.Text assignment is a no-op when the new text is
#        byte-identical to the text that was searched for.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute('{{ question.value | date("dd.MM.YYYY") }}', $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "___DATEFORMATPLACEHOLDER___"

    $rng2 = $d.Content
    $rng2.Find.Execute("___DATEFORMATPLACEHOLDER___", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $rng2.Text = '{{ question.value | date("dd.MM.YYYY") }}'
}

# ---------------------------------------------------------------------------
# 3) Signature table labels become template placeholders.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("Ort und Datum", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $rng3.Text = "{{signatureMetadata}}"
}

$rng4 = $d.Content
$found4 = $rng4.Find.Execute("Unterschrift", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found4) {
    $rng4.Text = "{{signatureTitle}}"
}
